$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking values so they are not
# auto-converted to actual numbers by Excel (column holds formatted text).
$textCells = @('D5', 'D6', 'D8', 'D9', 'D10', 'D11', 'D13', 'D14', 'D16', 'D18', 'D20', 'D22', 'D23', 'D26', 'D27', 'D28', 'D29', 'D30', 'D32', 'D33', 'D34', 'D35', 'D37', 'D39', 'D40', 'D44', 'D45', 'D47', 'D48', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '36.342.88'
$ws.Range('E2').Value = '  -1.48%  '

$ws.Range('D3').Value = '2.048.39'
$ws.Range('E3').Value = '  -1.90%  '

$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').Value = '245.48'
$ws.Range('E5').Value = '  -0.13%  '

$ws.Range('D6').Value = '0.661'
$ws.Range('E6').Value = '  +1.08%  '

$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').Value = '56.06'
$ws.Range('E8').Value = '  +0.44%  '

$ws.Range('D9').Value = '60.26'
$ws.Range('E9').Value = '  +0.67%  '

$ws.Range('D10').Value = '0.368'
$ws.Range('E10').Value = '  -0.34%  '

$ws.Range('D11').Value = '0.0746'
$ws.Range('E11').Value = '  -2.64%  '

$ws.Range('E12').Value = '  -3.59%  '

$ws.Range('D13').Value = '0.920'
$ws.Range('E13').Value = '  +3.22%  '

$ws.Range('D14').Value = '14.46'
$ws.Range('E14').Value = '  -4.41%  '

$ws.Range('D15').Value = '2.350.82'
$ws.Range('E15').Value = '  -1.81%  '

$ws.Range('D16').Value = '5.39'
$ws.Range('E16').Value = '  -2.69%  '

$ws.Range('D17').Value = '2.043.60'
$ws.Range('E17').Value = '  -3.97%  '

$ws.Range('D18').Value = '17.66'
$ws.Range('E18').Value = '  +1.10%  '

$ws.Range('D19').Value = '36.316.45'
$ws.Range('E19').Value = '  -1.48%  '

$ws.Range('D20').Value = '71.61'
$ws.Range('E20').Value = '  -2.21%  '

$ws.Range('D21').Value = '0.0₃0858'
$ws.Range('E21').Value = '  -3.02%  '

$ws.Range('D22').Value = '236.42'
$ws.Range('E22').Value = '  -0.58%  '

$ws.Range('D23').Value = '5.23'
$ws.Range('E23').Value = '  -4.93%  '

$ws.Range('E24').Value = '  +0.16%  '

$ws.Range('E25').Value = '  -2.56%  '

$ws.Range('D26').Value = '2.26'
$ws.Range('E26').Value = '  +4.10%  '

$ws.Range('D27').Value = '9.40'
$ws.Range('E27').Value = '  -5.06%  '

$ws.Range('D28').Value = '164.82'
$ws.Range('E28').Value = '  -2.28%  '

$ws.Range('D29').Value = '20.00'
$ws.Range('E29').Value = '  -3.21%  '

$ws.Range('D30').Value = '0.121'
$ws.Range('E30').Value = '  -1.76%  '

$ws.Range('E31').Value = '  -1.20%  '

$ws.Range('D32').Value = '4.98'
$ws.Range('E32').Value = '  -8.52%  '

$ws.Range('D33').Value = '0.0600'
$ws.Range('E33').Value = '  -2.04%  '

$ws.Range('D34').Value = '4.39'
$ws.Range('E34').Value = '  -7.02%  '

$ws.Range('D35').Value = '0.0895'
$ws.Range('E35').Value = '  +6.38%  '

$ws.Range('E36').Value = '  +0.09%  '

$ws.Range('D37').Value = '1.83'
$ws.Range('E37').Value = '  -0.40%  '

$ws.Range('E38').Value = '  -7.16%  '

$ws.Range('D39').Value = '5.06'
$ws.Range('E39').Value = '  +3.70%  '

$ws.Range('D40').Value = '1.21'
$ws.Range('E40').Value = '  -5.94%  '

$ws.Range('E41').Value = '  +1.44%  '

$ws.Range('E42').Value = '  -2.84%  '

$ws.Range('E43').Value = '  -5.44%  '

$ws.Range('D44').Value = '93.56'
$ws.Range('E44').Value = '  -3.67%  '

$ws.Range('D45').Value = '0.0907'
$ws.Range('E45').Value = '  -4.69%  '

$ws.Range('D46').Value = '1.398.99'
$ws.Range('E46').Value = '  +3.00%  '

$ws.Range('D47').Value = '15.87'
$ws.Range('E47').Value = '  -1.64%  '

$ws.Range('D48').Value = '7.48'
$ws.Range('E48').Value = '  +6.32%  '

$ws.Range('E49').Value = '  +1.73%  '

$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = '2.27'
$ws.Range('E50').Value = '  -7.55%  '

$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').Value = '46.24'
$ws.Range('E51').Value = '  +1.16%  '
